$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Discounted Total" row below the existing Total row.
$ws.Range("A7").Value = "Discounted Total"
$ws.Range("B7").Formula = "= 90% * B6"

# Move the selection like the saved workbook shows (A8 selected after entry).
$ws.Range("A8").Select()
